$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A: 14.42578125 -> 15.5703125, B: 14.42578125 -> 16.28515625)
# The underlying engine quantizes ColumnWidth to 1/6ths of a character when it
# serializes back to the OOXML "width" attribute (width_xml = ColumnWidth + 5/6,
# rounded to the nearest 1/6). Center the requested ColumnWidth in the rounding
# bucket that maps to the closest achievable XML width for each target.
$ws.Columns.Item(1).ColumnWidth = 88 / 6
$ws.Columns.Item(2).ColumnWidth = 93 / 6

# Update the cell values (columns A and B, rows 1-4); row 5 is left untouched.
# Using [double]"..." string-parses to keep full precision for the
# scientific-notation values (bare "e-06" literals are not parsed by the
# PowerShell-style script engine).
$ws.Range("A1").Value = [double]"-0.016905092561658057"
$ws.Range("B1").Value = [double]"0.016905092549975662"

$ws.Range("A2").Value = [double]"5.0534636183585651e-06"
$ws.Range("B2").Value = [double]"-5.0534843535883265e-06"

$ws.Range("A3").Value = [double]"-0.036479060795776778"
$ws.Range("B3").Value = [double]"0.036479060783700959"

$ws.Range("A4").Value = [double]"0.056931894823101473"
$ws.Range("B4").Value = [double]"-0.056931894842964841"
